$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Write new columns I:O with placement-vs-marketvalue data ---
$ws.Range("I1").Value = "Kader"
$ws.Range("J1").Value = "ø-Alter"
$ws.Range("K1").Value = "Legionäre"
$ws.Range("L1").Value = "Gesamtmarktwert"
$ws.Range("M1").Value = "ø-Marktwert"

$ws.Range("K2").Value = 644
$ws.Range("L2").Value = "23,5 Jahre"
$ws.Range("M2").Value = 320
$ws.Range("N2").Value = "2,58 Mrd. €"
$ws.Range("O2").Value = "4,01 Mio. €"

$ws.Range("I3").Value = "FC Bayern München"
$ws.Range("J3").Value = "FC Bayern München"
$ws.Range("K3").Value = 33
$ws.Range("L3").Value = "25,1"
$ws.Range("M3").Value = 17
$ws.Range("N3").Value = "595,40 Mio. €"
$ws.Range("O3").Value = "18,04 Mio. €"

$ws.Range("I4").Value = "Borussia Dortmund"
$ws.Range("J4").Value = "Borussia Dortmund"
$ws.Range("K4").Value = 31
$ws.Range("L4").Value = "24,2"
$ws.Range("M4").Value = 17
$ws.Range("N4").Value = "321,05 Mio. €"
$ws.Range("O4").Value = "10,36 Mio. €"

$ws.Range("I5").Value = "Bayer 04 Leverkusen"
$ws.Range("J5").Value = "Bayer 04 Leverkusen"
$ws.Range("K5").Value = 39
$ws.Range("L5").Value = "22,8"
$ws.Range("M5").Value = 21
$ws.Range("N5").Value = "244,58 Mio. €"
$ws.Range("O5").Value = "6,27 Mio. €"

$ws.Range("I6").Value = "FC Schalke 04"
$ws.Range("J6").Value = "FC Schalke 04"
$ws.Range("K6").Value = 39
$ws.Range("L6").Value = "23,7"
$ws.Range("M6").Value = 22
$ws.Range("N6").Value = "238,75 Mio. €"
$ws.Range("O6").Value = "6,12 Mio. €"

$ws.Range("I7").Value = "VfL Wolfsburg"
$ws.Range("J7").Value = "VfL Wolfsburg"
$ws.Range("K7").Value = 43
$ws.Range("L7").Value = "23,5"
$ws.Range("M7").Value = 21
$ws.Range("N7").Value = "225,35 Mio. €"
$ws.Range("O7").Value = "5,24 Mio. €"

$ws.Range("I8").Value = "Borussia Mönchengladbach"
$ws.Range("J8").Value = "Borussia Mönchengladbach"
$ws.Range("K8").Value = 33
$ws.Range("L8").Value = "23,4"
$ws.Range("M8").Value = 15
$ws.Range("N8").Value = "155,80 Mio. €"
$ws.Range("O8").Value = "4,72 Mio. €"

$ws.Range("I9").Value = "Hamburger SV"
$ws.Range("J9").Value = "Hamburger SV"
$ws.Range("K9").Value = 38
$ws.Range("L9").Value = "23,6"
$ws.Range("M9").Value = 20
$ws.Range("N9").Value = "94,98 Mio. €"
$ws.Range("O9").Value = "2,50 Mio. €"

$ws.Range("I10").Value = "1.FSV Mainz 05"
$ws.Range("J10").Value = "1.FSV Mainz 05"
$ws.Range("K10").Value = 39
$ws.Range("L10").Value = "23,3"
$ws.Range("M10").Value = 22
$ws.Range("N10").Value = "85,70 Mio. €"
$ws.Range("O10").Value = "2,20 Mio. €"

$ws.Range("I11").Value = "1.FC Köln"
$ws.Range("J11").Value = "1.FC Köln"
$ws.Range("K11").Value = 30
$ws.Range("L11").Value = "24,4"
$ws.Range("M11").Value = 14
$ws.Range("N11").Value = "77,40 Mio. €"
$ws.Range("O11").Value = "2,58 Mio. €"

$ws.Range("I12").Value = "TSG 1899 Hoffenheim"
$ws.Range("J12").Value = "TSG 1899 Hoffenheim"
$ws.Range("K12").Value = 33
$ws.Range("L12").Value = "23,5"
$ws.Range("M12").Value = 15
$ws.Range("N12").Value = "76,55 Mio. €"
$ws.Range("O12").Value = "2,32 Mio. €"

$ws.Range("I13").Value = "Hertha BSC"
$ws.Range("J13").Value = "Hertha BSC"
$ws.Range("K13").Value = 33
$ws.Range("L13").Value = "24,0"
$ws.Range("M13").Value = 14
$ws.Range("N13").Value = "75,98 Mio. €"
$ws.Range("O13").Value = "2,30 Mio. €"

$ws.Range("I14").Value = "SV Werder Bremen"
$ws.Range("J14").Value = "SV Werder Bremen"
$ws.Range("K14").Value = 41
$ws.Range("L14").Value = "24,5"
$ws.Range("M14").Value = 20
$ws.Range("N14").Value = "75,08 Mio. €"
$ws.Range("O14").Value = "1,83 Mio. €"

$ws.Range("I15").Value = "RasenBallsport Leipzig"
$ws.Range("J15").Value = "RasenBallsport Leipzig  "
$ws.Range("K15").Value = 33
$ws.Range("L15").Value = "22,6"
$ws.Range("M15").Value = 15
$ws.Range("N15").Value = "69,18 Mio. €"
$ws.Range("O15").Value = "2,10 Mio. €"

$ws.Range("I16").Value = "FC Augsburg"
$ws.Range("J16").Value = "FC Augsburg"
$ws.Range("K16").Value = 37
$ws.Range("L16").Value = "25,3"
$ws.Range("M16").Value = 21
$ws.Range("N16").Value = "69,00 Mio. €"
$ws.Range("O16").Value = "1,86 Mio. €"

$ws.Range("I17").Value = "Eintracht Frankfurt"
$ws.Range("J17").Value = "Eintracht Frankfurt"
$ws.Range("K17").Value = 35
$ws.Range("L17").Value = "23,1"
$ws.Range("M17").Value = 21
$ws.Range("N17").Value = "66,10 Mio. €"
$ws.Range("O17").Value = "1,89 Mio. €"

$ws.Range("I18").Value = "SC Freiburg"
$ws.Range("J18").Value = "SC Freiburg  Deutscher Zweitligameister 15/16"
$ws.Range("K18").Value = 36
$ws.Range("L18").Value = "24,0"
$ws.Range("M18").Value = 13
$ws.Range("N18").Value = "40,25 Mio. €"
$ws.Range("O18").Value = "1,12 Mio. €"

$ws.Range("I19").Value = "FC Ingolstadt 04"
$ws.Range("J19").Value = "FC Ingolstadt 04"
$ws.Range("K19").Value = 32
$ws.Range("L19").Value = "23,8"
$ws.Range("M19").Value = 15
$ws.Range("N19").Value = "36,68 Mio. €"
$ws.Range("O19").Value = "1,15 Mio. €"

$ws.Range("I20").Value = "SV Darmstadt 98"
$ws.Range("J20").Value = "SV Darmstadt 98"
$ws.Range("K20").Value = 39
$ws.Range("L20").Value = "25,7"
$ws.Range("M20").Value = 17
$ws.Range("N20").Value = "32,83 Mio. €"
$ws.Range("O20").Value = "842 Tsd. €"

# --- Column widths for new columns (closest match to bestFit autofit widths) ---
$ws.Columns.Item(9).ColumnWidth = 24.5
$ws.Columns.Item(10).ColumnWidth = 66.0
$ws.Columns.Item(11).ColumnWidth = 8.8333
$ws.Columns.Item(12).ColumnWidth = 16.3333
$ws.Columns.Item(13).ColumnWidth = 11.3333
$ws.Columns.Item(14).ColumnWidth = 11.6667

# --- Selection used by the author while reviewing the new data ---
$ws.Range("O3:O20").Select()

# --- Page setup matching the printed layout ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

